$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in T2 (129731 -> 130416)
$ws.Range("T2").Value = 130416

# Move selection to T3 (activeCell/sqref T2 -> T3)
$ws.Range("T3").Select()
